# Generate Report for Handoff
# The "af843c2a-ac76-495f-a669-548f82605fb1.md" file is now ready for
# handoff: update its Status/Latest-Handoff fields on the Overview sheet
# and on each per-locale sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for af843c2a-...md (row 3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-37-19 22:37:02"

# --- zh-cn sheet: row for af843c2a-...md (row 3) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-19 22:36:59"

# --- de-de sheet: row for af843c2a-...md (row 3) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-19 22:37:02"
